$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order / rewrite the header row (row 1) and add the new "Urgenta" column (E) ---
# Copy the existing header's style (bold font, border, centered alignment) onto the
# new E1 header cell before writing its text, so it matches the other headers (s="1").
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Value = "Nume"
$ws.Range("B1").Value = "Prenume"
$ws.Range("C1").Value = "Data Tiparire"
$ws.Range("D1").Value = "Perioada Internarii"
$ws.Range("E1").Value = "Urgenta"

# --- Row 2 (existing record, columns re-shuffled + new Urgenta value) ---
$ws.Range("A2").Value = "AFTINIE `n"
$ws.Range("B2").Value = "GRIGORITA `n"
$ws.Range("C2").Value = "04/08/2022 "
$ws.Range("D2").Value = "26/11/2019 08:04 - 02/12/2019 10:05 (6 zile) "
$ws.Range("E2").Value = "NU `n"

# --- Row 3 (new record) ---
$ws.Range("A3").Value = "ALIMAN `n"
$ws.Range("B3").Value = "MOISE `n"
$ws.Range("C3").Value = "04/08/2022 "
$ws.Range("D3").Value = "05/12/2019 17:03 - 13/12/2019 12:20 (8 zile) "
$ws.Range("E3").Value = "DA `n"

# --- Row 4 (new record) ---
$ws.Range("A4").Value = "AVRAM `n"
$ws.Range("B4").Value = "IULICA `n"
$ws.Range("C4").Value = "04/08/2022 "
$ws.Range("D4").Value = "18/11/2019 09:20 - 22/11/2019 13:34 (4 zile) "
$ws.Range("E4").Value = "NU `n"

# --- Row 5 (new record) ---
$ws.Range("A5").Value = "BALASA `n"
$ws.Range("B5").Value = "IOANA `n"
$ws.Range("C5").Value = "04/08/2022 "
$ws.Range("D5").Value = "25/11/2019 09:19 - 29/11/2019 12:44 (4 zile) "
$ws.Range("E5").Value = "NU `n"

# --- Row 6 (trailing blank row present in the sheet's used range) ---
# Touch the cells (without actually assigning content) so the engine registers
# them as present-but-empty and the sheet dimension grows to A1:E6, matching
# the source file's empty <c .../> placeholders.
$ws.Range("A6:E6").Font.Bold = $false
